$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.883.70'
$ws.Range("E2").Value = '  -5.05%  '
$ws.Range("D3").Value = '2.493.51'
$ws.Range("E3").Value = '  -3.24%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '536.34'
$ws.Range("E5").Value = '  -2.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.94'
$ws.Range("E6").Value = '  -7.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.35%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.573'
$ws.Range("E8").Value = '  -3.92%  '
$ws.Range("D9").Value = '2.513.75'
$ws.Range("E9").Value = '  -2.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.100'
$ws.Range("E10").Value = '  -3.95%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.55'
$ws.Range("E12").Value = '  +0.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.352'
$ws.Range("E13").Value = '  -4.27%  '
$ws.Range("D14").Value = '2.927.53'
$ws.Range("E14").Value = '  -3.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.91'
$ws.Range("E15").Value = '  -7.14%  '
$ws.Range("D16").Value = '58.815.22'
$ws.Range("E16").Value = '  -4.98%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000139'
$ws.Range("E17").Value = '  -4.35%  '
$ws.Range("D18").Value = '2.505.77'
$ws.Range("E18").Value = '  -2.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.31'
$ws.Range("E19").Value = '  -2.72%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.28'
$ws.Range("E20").Value = '  -6.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '322.70'
$ws.Range("E21").Value = '  -4.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.74'
$ws.Range("E23").Value = '  -5.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.85'
$ws.Range("E24").Value = '  -4.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.440'
$ws.Range("E25").Value = '  -11.03%  '
$ws.Range("E26").Value = '  -3.68%  '
$ws.Range("D27").Value = '2.613.63'
$ws.Range("E27").Value = '  -3.09%  '
$ws.Range("E28").Value = '  -0.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.75'
$ws.Range("E29").Value = '  -5.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.85'
$ws.Range("E30").Value = '  -6.98%  '
$ws.Range("D31").Value = '0.0₃0779'
$ws.Range("E31").Value = '  -7.53%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.25'
$ws.Range("E32").Value = '  -7.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.78'
$ws.Range("E33").Value = '  -6.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.996'
$ws.Range("E34").Value = '  -0.26%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '155.93'
$ws.Range("E35").Value = '  -4.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.43'
$ws.Range("E36").Value = '  -1.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.50'
$ws.Range("E37").Value = '  -3.94%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.40'
$ws.Range("E38").Value = '  -9.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.61'
$ws.Range("E39").Value = '  -10.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.90'
$ws.Range("E40").Value = '  -2.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '309.17'
$ws.Range("E41").Value = '  -6.57%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '36.76'
$ws.Range("E42").Value = '  -2.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.67'
$ws.Range("E43").Value = '  -7.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.811'
$ws.Range("E44").Value = '  -11.55%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.990'
$ws.Range("E45").Value = '  -0.74%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.594'
$ws.Range("E46").Value = '  -2.19%  '
$ws.Range("E47").Value = '  -1.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.06'
$ws.Range("E48").Value = '  +1.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0928'
$ws.Range("E49").Value = '  -4.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.53'
$ws.Range("E50").Value = '  -5.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0516'
$ws.Range("E51").Value = '  -6.10%  '
